$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the last-revised date shown in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "FPIEBP" sheet: update hard coal's balancing priorities (row 3) ---
$wsData = $wb.Worksheets.Item("FPIEBP")
$wsData.Range("B3").Value = 1
$wsData.Range("C3").Value = 3
$wsData.Range("D3").Value = 2

# Move the active selection on the FPIEBP sheet to E3 (matches the saved cursor position)
$wsData.Activate()
$wsData.Range("E3").Select() | Out-Null
